$wb = $excel.ActiveWorkbook

$wsTestData = $wb.Worksheets.Item("TestData")
$wsUS       = $wb.Worksheets.Item("US")
$wsDemo     = $wb.Worksheets.Item("Demo")
$wsBackup   = $wb.Worksheets.Item("backup")
$wsData     = $wb.Worksheets.Item("Data")

# --- 1. Archive the current "TestData" row 2 (EU/DE case) into the "Data"
#        sheet (row 5) before it gets overwritten with the new test case.
$wsTestData.Range("A2:O2").Copy($wsData.Range("A5:O5"))

# --- 2. Duplicate some of the Data sheet's existing rows (2,3,4) further
#        down (rows 9/15, 10/16 and 11) as additional D365 API test cases.
$wsData.Range("A2:O2").Copy($wsData.Range("A9:O9"))
$wsData.Range("A3:O3").Copy($wsData.Range("A10:O10"))
$wsData.Range("A4:O4").Copy($wsData.Range("A11:O11"))
$wsData.Range("A2:O2").Copy($wsData.Range("A15:O15"))
$wsData.Range("A3:O3").Copy($wsData.Range("A16:O16"))

$wsData.Rows.Item(9).RowHeight = 60
$wsData.Rows.Item(10).RowHeight = 60
$wsData.Rows.Item(11).RowHeight = 60
$wsData.Rows.Item(15).RowHeight = 60
$wsData.Rows.Item(16).RowHeight = 60

# --- 3. Overwrite "TestData" row 2 with the new "US region - Paypal"
#        D365 API test case.
$wsTestData.Range("A2").Value = "Place one FG & subscription order using Paypal as a payment method from Breville and validate the order details in My Breville drop down - US/EN"
$wsTestData.Range("B2").Value = "us"
$wsTestData.Range("C2").Value = "NA"
$wsTestData.Range("D2").Value = "Mahesh"
$wsTestData.Range("E2").Value = "N"
$wsTestData.Range("F2").Value = "111 West Adams Street"
$wsTestData.Range("G2").Value = "Chicago"
$wsTestData.Range("H2").Value = 60603
$wsTestData.Range("I2").Value = "Illinois"
$wsTestData.Range("J2").Value = 7022624000
$wsTestData.Range("K2").Value = "NewUser"
$wsTestData.Range("L2").Value = "NA"
$wsTestData.Range("M2").Value = "test12345"
$wsTestData.Range("N2").Value = "475 Yonge St"
$wsTestData.Range("O2").Value = "Creditcard"

$wsTestData.Range("A2").WrapText = $true
$wsTestData.Rows.Item(2).RowHeight = 60

# --- 4. Update the selections recorded for each sheet view.
$wsTestData.Activate()
$wsTestData.Range("P2").Select()

$wsUS.Activate()
$wsUS.Range("O2").Select()

$wsDemo.Activate()
$wsDemo.Range("K2").Select()

$wsBackup.Activate()
$wsBackup.Range("O46").Select()

$wsData.Activate()
$wsData.Range("O16").Select()

$wsTestData.Activate()
